$wb = $excel.ActiveWorkbook

# --- Update selections on existing sheets (2d, 3d) ---
$ws2 = $wb.Worksheets.Item("2d")
[void]$ws2.Range("C8").Select()

$ws3 = $wb.Worksheets.Item("3d")
[void]$ws3.Range("B15").Select()

# --- Add the new "missing_values" sheet after the last sheet ("5d") ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "missing_values"

# --- Fill in the data (age / sex\time table with one missing combination) ---
$ws.Range("A1").Value = "age"
$ws.Range("B1").Value = "sex\time"
$ws.Range("C1").Value = 2007
$ws.Range("D1").Value = 2010
$ws.Range("E1").Value = 2013

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "F"
$ws.Range("C2").Value = 3722
$ws.Range("D2").Value = 3395
$ws.Range("E2").Value = 3347

$ws.Range("A3").Value = 0
$ws.Range("B3").Value = "H"
$ws.Range("C3").Value = 338
$ws.Range("D3").Value = 316
$ws.Range("E3").Value = 323

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "F"
$ws.Range("C4").Value = 2878
$ws.Range("D4").Value = 2791
$ws.Range("E4").Value = 2822

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "F"
$ws.Range("C5").Value = 4073
$ws.Range("D5").Value = 4161
$ws.Range("E5").Value = 4429

$ws.Range("A6").Value = 2
$ws.Range("B6").Value = "H"
$ws.Range("C6").Value = 1561
$ws.Range("D6").Value = 1463
$ws.Range("E6").Value = 1467

$ws.Range("A7").Value = 3
$ws.Range("B7").Value = "F"
$ws.Range("C7").Value = 3507
$ws.Range("D7").Value = 3741
$ws.Range("E7").Value = 3366

$ws.Range("A8").Value = 3
$ws.Range("B8").Value = "H"
$ws.Range("C8").Value = 2052
$ws.Range("D8").Value = 2052
$ws.Range("E8").Value = 2118

$ws.Range("A9").Value = 4
$ws.Range("B9").Value = "H"
$ws.Range("C9").Value = 3785
$ws.Range("D9").Value = 3508
$ws.Range("E9").Value = 3172

# --- Select the cell that was active on the new sheet and make it the active/selected tab ---
[void]$ws.Range("G24").Select()
